$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New time-log entry for row 72.
# Set D72 (Interruption minutes) first so the shared formula in E72
# recalculates against the final value instead of a stale default.
$ws.Range("D72").Value = 20
$ws.Range("A72").Value = 41916
$ws.Range("B72").Value = 0.95833333333333337
$ws.Range("C72").Value = 1.01875
$ws.Range("F72").Value = "Coding"

# Move the active selection the way the author's session ended up.
$ws.Range("C73").Select()
